$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.618.40"
$ws.Range("D3").Value = "1.644.17"
$ws.Range("E3").Value = "  +0.79%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").Value = "215.89"
$ws.Range("E5").Value = "  +1.63%  "
$ws.Range("E6").Value = "  +0.99%  "
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("D9").Value = "0.0627"
$ws.Range("E9").Value = "  +0.91%  "
$ws.Range("D10").Value = "19.21"
$ws.Range("E10").Value = "  +0.39%  "
$ws.Range("E11").Value = "  -0.02%  "
$ws.Range("D12").Value = "1.872.29"
$ws.Range("E12").Value = "  +0.69%  "
$ws.Range("D13").Value = "4.22"
$ws.Range("E13").Value = "  +3.44%  "
$ws.Range("D14").Value = "1.621.51"
$ws.Range("E14").Value = "  -1.05%  "
$ws.Range("E15").Value = "  +1.64%  "
$ws.Range("D16").Value = "65.89"
$ws.Range("E16").Value = "  +4.42%  "
$ws.Range("D17").Value = "26.661.80"
$ws.Range("E17").Value = "  +0.06%  "
$ws.Range("E18").Value = "  +1.47%  "
$ws.Range("D19").Value = "218.36"
$ws.Range("E19").Value = "  +0.46%  "
$ws.Range("E20").Value = "  +0.30%  "
$ws.Range("E21").Value = "  +2.61%  "
$ws.Range("D22").Value = "6.32"
$ws.Range("E22").Value = "  +2.48%  "
$ws.Range("D23").Value = "9.55"
$ws.Range("E23").Value = "  +2.06%  "
$ws.Range("D24").Value = "2.13"
$ws.Range("E24").Value = "  +8.84%  "
$ws.Range("D25").Value = "146.42"
$ws.Range("E25").Value = "  -1.07%  "
$ws.Range("E26").Value = "  +0.38%  "
$ws.Range("E27").Value = "  -0.27%  "
$ws.Range("D28").Value = "7.13"
$ws.Range("E28").Value = "  +4.33%  "
$ws.Range("E29").Value = "  +2.45%  "
$ws.Range("D30").Value = "0.0519"
$ws.Range("E30").Value = "  +2.96%  "
$ws.Range("E31").Value = "  +0.55%  "
$ws.Range("E33").Value = "  +2.88%  "
$ws.Range("D34").Value = "1.275.11"
$ws.Range("E34").Value = "  +5.41%  "
$ws.Range("E35").Value = "  +2.64%  "
$ws.Range("D36").Value = "0.0182"
$ws.Range("E36").Value = "  +6.00%  "
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("D38").Value = "0.531"
$ws.Range("E38").Value = "  +6.55%  "
$ws.Range("E39").Value = "  +2.57%  "
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.24%  "
$ws.Range("E41").Value = "  +2.18%  "
$ws.Range("E42").Value = "  -1.50%  "
$ws.Range("D43").Value = "5.48"
$ws.Range("E43").Value = "  +1.46%  "
$ws.Range("D44").Value = "1.783.71"
$ws.Range("E44").Value = "  +0.41%  "
$ws.Range("D45").Value = "93.34"
$ws.Range("E45").Value = "  +0.41%  "
$ws.Range("D46").Value = "59.72"
$ws.Range("E46").Value = "  +9.36%  "
$ws.Range("E47").Value = "  +3.90%  "
$ws.Range("E48").Value = "  +0.96%  "
$ws.Range("D49").Value = "7.79"
$ws.Range("E49").Value = "  +2.91%  "
$ws.Range("E51").Value = "  -0.56%  "
